$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Set all Runmode cells (C2:C7) to "Y" - running all suites
$ws.Range("C2:C7").Value = "Y"

# Update the active selection to C8
$ws.Range("C8").Select()
